$wb = $excel.ActiveWorkbook

# Sheet "建物" (building) - 2nd worksheet: column I (property_category) rows 2-9
# Currently holds "land", should be "building"
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 9; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

# Sheet "汽車" (car) - 3rd worksheet: column H (property_category) row 2
# Currently holds "land", should be "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Cells.Item(2, 8).Value = "car"
